$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.118.52"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "'3.079.14"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'543.70"
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("D6").Value = "'137.25"
$ws.Range("E6").Value = "  -6.42%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'3.070.84"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "'6.54"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'34.72"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "'0.0000217"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'3.568.03"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "'63.096.00"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'3.071.53"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "'492.42"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "'0.698"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("D24").Value = "'77.51"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'12.22"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'8.32"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "'1.93"
$ws.Range("E30").Value = "  -8.57%  "
$ws.Range("D31").Value = "'26.33"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("D34").Value = "'59.82"
$ws.Range("E34").Value = "  +13.16%  "
$ws.Range("D35").Value = "'522.56"
$ws.Range("E35").Value = "  -6.92%  "
$ws.Range("D36").Value = "'5.92"
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").Value = "'5.16"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").Value = "'0.0401"
$ws.Range("E38").Value = "  -6.90%  "
$ws.Range("D39").Value = "'0.0790"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "'3.041.27"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "'8.08"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").Value = "'0.255"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "  -6.93%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'121.85"
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("D48").Value = "'24.23"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "'0.107"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'0.0₃0506"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").Value = "'2.32"
$ws.Range("E51").Value = "  +38.17%  "
